# RefinementLevels.xlsx update
# - Adds a new "Level 1" derived column (H8 = G8/15)
# - Updates the refinement-ratio inputs for Levels 4-7 (columns B and J)
#   which drives the recalculation of the dependent formula columns
#   C:G and K:O for rows 11-14
# - Documents the South/North boundary-forcing notes in column L
# - Updates the active selection to reflect where the user was working

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New formula cell in row 8 (Level 1 block): seconds per 15 -> H8
$ws.Range("H8").Formula = "=G8/15"

# Updated refinement-ratio counts (left block, columns A-G, "Lon" side)
$ws.Range("B11").Value = 6   # Level 4
$ws.Range("B12").Value = 7   # Level 5
$ws.Range("B13").Value = 11  # Level 6
$ws.Range("B14").Value = 2   # Level 7

# Updated refinement-ratio counts (right block, columns I-O, "Lat" side)
$ws.Range("J11").Value = 6   # Level 4
$ws.Range("J13").Value = 2   # Level 6
$ws.Range("J14").Value = 5   # Level 7

# New annotation cells (order matters so the shared-string table is
# built in the same sequence as the target workbook)
$ws.Range("L22").Value = "forcing the depth, copy the boundaries into the domain then set the depth"
$ws.Range("L21").Value = "set the height and depth based on the tide"

# Reflect the cell the user was last working in
$ws.Range("B13").Select()
